$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.864.33"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.758.02"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'327.70"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.4689"
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("D8").Value = "'0.3502"
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("D9").Value = "'42.26"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "'0.07356"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").Value = "'1.082"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'1.0000"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'20.51"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").Value = "'5.982"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "'7.149"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "1.753.56"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "'92.06"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'0.06409"
$ws.Range("D21").Value = "'16.77"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").Value = "'5.750"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "27.883.57"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").Value = "'2.150"
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").Value = "'161.71"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "1.956.95"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "'2.152"
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").Value = "'122.69"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D31").Value = "'1.071"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'0.09345"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").Value = "'3.652"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'5.545"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "'0.02267"
$ws.Range("D36").Value = "'11.62"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Value = "'4.897"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "'0.6136"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("D41").Value = "'1.177"
$ws.Range("D42").Value = "'7.773"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "'1.353"
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("D44").Value = "'13.02"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").Value = "'3.733"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "'0.5782"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "'122.59"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'1.922"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").Value = "'1.120"
$ws.Range("D51").Value = "'72.07"
$ws.Range("E51").Value = "  -0.78%  "
